# Update the "Förändrad" date column (C2:C6) from 2023-10-05 (45204) to 2023-10-08 (45207)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..6) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45207
    }
}
